$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 240
$ws.Range("F3").Value = 1397
$ws.Range("F5").Value = 890
$ws.Range("F7").Value = 1234
$ws.Range("F8").Value = 1557
$ws.Range("F9").Value = 158
$ws.Range("F11").Value = 2258
$ws.Range("F12").Value = 447
$ws.Range("F13").Value = 114
$ws.Range("F15").Value = 28
$ws.Range("F16").Value = 88
$ws.Range("F17").Value = 81
$ws.Range("F18").Value = 6137
$ws.Range("F20").Value = 6004
$ws.Range("F21").Value = 9979
$ws.Range("F23").Value = 172
$ws.Range("F24").Value = 182
$ws.Range("F25").Value = 275
$ws.Range("F26").Value = 495
$ws.Range("F27").Value = 165
$ws.Range("F28").Value = 146
$ws.Range("F29").Value = 4385
$ws.Range("F30").Value = 382

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 789
$ws.Range("F3").Value = 609

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 789
$ws.Range("F3").Value = 609
$ws.Range("F4").Value = 240
$ws.Range("F5").Value = 1397
$ws.Range("F8").Value = 890
$ws.Range("F10").Value = 1234
$ws.Range("F12").Value = 1557
$ws.Range("F14").Value = 158
$ws.Range("F15").Value = 2258
$ws.Range("F17").Value = 447
$ws.Range("F18").Value = 114
$ws.Range("F20").Value = 28
$ws.Range("F22").Value = 88
$ws.Range("F23").Value = 81
$ws.Range("F24").Value = 6137
$ws.Range("F26").Value = 6004
$ws.Range("F27").Value = 9979
$ws.Range("F30").Value = 172
$ws.Range("F31").Value = 182
$ws.Range("F32").Value = 275
$ws.Range("F34").Value = 495
$ws.Range("F38").Value = 165
$ws.Range("F39").Value = 146
$ws.Range("F40").Value = 4385
$ws.Range("F46").Value = 382
